## "day15 内置和开发规范/代码和作业答案/课堂代码/db.xlsx"
## Add more course material: new student rows/columns for the little
## "user db" demo sheet, plus the default-font/selection tidy-up that
## came along with regenerating this sheet's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- replace the 3 existing rows with the refreshed roster, and add a
#     4th column (short token) for each student --------------------------
$ws.Range("A1").Value = "Lili"
$ws.Range("B1").Value = "24034476e99d1aa26614c9e7902ae2700316c6709aeaedb24e79231968521b7d"
$ws.Range("C1").Value = "2022-02-25 14:29:10"
$ws.Range("D1").Value = "A4nf3I8M"

$ws.Range("A2").Value = "Achuan"
$ws.Range("B2").Value = "91dfd8e0b70f8008238d780a57a64d17b5758031dd4b69ec9122644a3a45ea3c"
$ws.Range("C2").Value = "2022-02-25 14:29:14"
$ws.Range("D2").Value = "R9HrMT2z"

$ws.Range("A3").Value = "Achuan-2"
$ws.Range("B3").Value = "d01ea24f3c2e92c641b8e142cd3636c694690d7d655b0c3295f358fc9ac81ffc"
$ws.Range("C3").Value = "2022-02-25 14:31:18"
$ws.Range("D3").Value = "v3GjYw81"

# --- default workbook font: drop the locale (Chinese) default font in
#     favour of the standard Calibri used by the refreshed file ---------
try {
    $wb.Styles.Item(1).Font.Name = "Calibri"
} catch {
}

# --- reset the view back to the top-left cell (no stale selection/tab
#     leftover from the previous editing session) ------------------------
[void]$ws.Range("A1").Select()
